# rev 2 changes, see revision_log
#
# Summary of the edit applied to the SMD component-placement (CPL) sheet:
#   1. Re-measured Mid X / Mid Y coordinates for C3, C4, C5, C6, C8, C9, C10.
#   2. A new component "Q2" is inserted into the placement list (between Q1
#      and the R-series), pushing every row from the old R1 down by one.
#   3. R10's coordinates were additionally updated once it moved into its
#      new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update re-measured Mid X / Mid Y for a handful of existing parts ---
# (Designator / Layer / Rotation stay the same; only B (Mid X) and C (Mid Y) change.)

$ws.Range("B4").Value = -14.2       # C3  Mid X
$ws.Range("C4").Value = 70.125      # C3  Mid Y

$ws.Range("B5").Value = -23.8       # C4  Mid X
$ws.Range("C5").Value = 66.325      # C4  Mid Y

$ws.Range("B6").Value = -14.2       # C5  Mid X
$ws.Range("C6").Value = 48.575      # C5  Mid Y

$ws.Range("B7").Value = -23.8       # C6  Mid X
$ws.Range("C7").Value = 48.025      # C6  Mid Y

$ws.Range("B9").Value = -14.45      # C8  Mid X
$ws.Range("C9").Value = 39.375      # C8  Mid Y

$ws.Range("B10").Value = -24.05     # C9  Mid X
$ws.Range("C10").Value = 39.525     # C9  Mid Y

$ws.Range("B11").Value = -11.938    # C10 Mid X (Mid Y unchanged)

# --- 2. Insert a new row for "Q2" right after "Q1" (row 14), pushing the ---
#        rest of the table (old rows 15-42, R1..U3) down by one row.
$ws.Rows(15).Insert()

$ws.Range("A15").Value = "Q2"
$ws.Range("B15").Value = -24.75
$ws.Range("C15").Value = 53.975
$ws.Range("D15").Value = "top"
$ws.Range("E15").Value = 270.0

# --- 3. R10 (now at row 25 after the insert above) gets new coordinates ---
$ws.Range("B25").Value = -24.7125   # R10 Mid X
$ws.Range("C25").Value = 45.35      # R10 Mid Y
